# Applies the "Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta" update:
# refreshes the EC (Estado de Cuenta) table on Hoja1 rows 16-47 (columns B:G) with the
# reorganized / extended employee-period dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row layout: B=Tipo Doc, C=N Doc, D=Nombre, E=Periodo Mora, F=Valor Mora, G=Salario Basico
$data = @(
    @(16, 'CC', '73508103',   'YORBIS ANTONIO ROSADO MENDOZA',      '2503', 79075,  1976894),
    @(17, 'CC', '73508103',   'YORBIS ANTONIO ROSADO MENDOZA',      '2502', 79075,  1976894),
    @(18, 'CC', '73508103',   'YORBIS ANTONIO ROSADO MENDOZA',      '2501', 79075,  1976894),
    @(19, 'CC', '73508103',   'YORBIS ANTONIO ROSADO MENDOZA',      '2412', 79075,  1976894),
    @(20, 'CC', '73508103',   'YORBIS ANTONIO ROSADO MENDOZA',      '2411', 79075,  1976894),
    @(21, 'CC', '73508103',   'YORBIS ANTONIO ROSADO MENDOZA',      '2410', 79075,  1976894),
    @(22, 'CC', '73508103',   'YORBIS ANTONIO ROSADO MENDOZA',      '2409', 79075,  1976894),
    @(23, 'CC', '73508103',   'YORBIS ANTONIO ROSADO MENDOZA',      '2408', 79075,  1976894),
    @(24, 'CC', '73508103',   'YORBIS ANTONIO ROSADO MENDOZA',      '2407', 79075,  1976894),
    @(25, 'CC', '73508103',   'YORBIS ANTONIO ROSADO MENDOZA',      '2406', 79075,  1976894),
    @(26, 'CC', '1065625378', 'GENDRIS ZULEIMA OROZCO RODRIGUEZ',   '2408', 120000, 3000000),
    @(27, 'CC', '73508103',   'YORBIS ANTONIO ROSADO MENDOZA',      '2409', 79075,  1976894),
    @(28, 'CC', '19596360',   'GERMAN MARTINEZ CALDERON',           '2409', 180000, 4500000),
    @(29, 'CC', '1065625378', 'GENDRIS ZULEIMA OROZCO RODRIGUEZ',   '2409', 120000, 3000000),
    @(30, 'CC', '73508103',   'YORBIS ANTONIO ROSADO MENDOZA',      '2410', 79075,  1976894),
    @(31, 'CC', '19596360',   'GERMAN MARTINEZ CALDERON',           '2410', 180000, 4500000),
    @(32, 'CC', '1065625378', 'GENDRIS ZULEIMA OROZCO RODRIGUEZ',   '2410', 120000, 3000000),
    @(33, 'CC', '73508103',   'YORBIS ANTONIO ROSADO MENDOZA',      '2411', 79075,  1976894),
    @(34, 'CC', '19596360',   'GERMAN MARTINEZ CALDERON',           '2411', 180000, 4500000),
    @(35, 'CC', '1065625378', 'GENDRIS ZULEIMA OROZCO RODRIGUEZ',   '2411', 120000, 3000000),
    @(36, 'CC', '73508103',   'YORBIS ANTONIO ROSADO MENDOZA',      '2412', 79075,  1976894),
    @(37, 'CC', '19596360',   'GERMAN MARTINEZ CALDERON',           '2412', 180000, 4500000),
    @(38, 'CC', '1065625378', 'GENDRIS ZULEIMA OROZCO RODRIGUEZ',   '2412', 120000, 3000000),
    @(39, 'CC', '73508103',   'YORBIS ANTONIO ROSADO MENDOZA',      '2501', 79075,  1976894),
    @(40, 'CC', '19596360',   'GERMAN MARTINEZ CALDERON',           '2501', 180000, 4500000),
    @(41, 'CC', '1065625378', 'GENDRIS ZULEIMA OROZCO RODRIGUEZ',   '2501', 120000, 3000000),
    @(42, 'CC', '73508103',   'YORBIS ANTONIO ROSADO MENDOZA',      '2502', 79075,  1976894),
    @(43, 'CC', '19596360',   'GERMAN MARTINEZ CALDERON',           '2502', 180000, 4500000),
    @(44, 'CC', '1065625378', 'GENDRIS ZULEIMA OROZCO RODRIGUEZ',   '2502', 120000, 3000000),
    @(45, 'CC', '73508103',   'YORBIS ANTONIO ROSADO MENDOZA',      '2503', 47445,  1976894),
    @(46, 'CC', '19596360',   'GERMAN MARTINEZ CALDERON',           '2503', 108000, 4500000),
    @(47, 'CC', '1065625378', 'GENDRIS ZULEIMA OROZCO RODRIGUEZ',   '2503', 72000,  3000000)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}
